# "fixed some errors after first run app"
# Correct the Stop loss and Position values that were wrong on the first run,
# and leave the selection on the Position cell that was just fixed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stop loss (C2): 72000 -> 70100
$ws.Range("C2").Value = 70100

# Position (E2): 0.01 -> 0.001
$ws.Range("E2").Value = 0.001

# Reflect the sheet's default/standard column width nudging that happened
# alongside these edits (best-effort; some runtimes may not persist this).
try { $ws.StandardWidth = 11.66015625 } catch {}

# Active cell/selection ends up on E2 after the fix
$ws.Range("E2").Select()
